$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet Hoja1 -> persona
$ws.Name = "persona"

# Add the new "frank" persona row (row 2)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "frank"
$ws.Range("C2").Value = 21
$ws.Range("D2").Value = 157799
$ws.Range("E2").Value = "frankrcuetia@gmail.com"
$ws.Range("F2").Value = "24/02/2005"

# correo (E2) becomes a mailto hyperlink - this also creates the
# Hyperlink cell style (underlined font, theme color 10) used for E2
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:frankrcuetia@gmail.com") | Out-Null

# nacimiento (F2) gets a date number format (builtin format id 14)
$ws.Range("F2").NumberFormat = "mm-dd-yy"

# Update selection to match the saved view state
$ws.Range("F8").Select() | Out-Null

# Page orientation was explicitly set to portrait
$ws.PageSetup.Orientation = 1

Write-Host "persona row added"
